$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.35
$ws.Range("H2").Value = 2.42
$ws.Range("O2").Value = 1.34
$ws.Range("V2").Value = 1.69
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 10.5
$ws.Range("AN2").Value = 36
$ws.Range("H3").Value = 4.5
$ws.Range("L3").Value = 1.27
$ws.Range("N3").Value = 4.3
$ws.Range("P3").Value = 2.14
$ws.Range("Q3").Value = 1.65
$ws.Range("R3").Value = 1.47
$ws.Range("AK3").Value = 1000
$ws.Range("F4").Value = 1.82
$ws.Range("G4").Value = 1.99
$ws.Range("H4").Value = 4.1
$ws.Range("J4").Value = 3.7
$ws.Range("T4").Value = 1.71
$ws.Range("W4").Value = 2
$ws.Range("Y4").Value = 22
$ws.Range("Z4").Value = 42
$ws.Range("AA4").Value = 120
$ws.Range("AB4").Value = 12
$ws.Range("AF4").Value = 15
$ws.Range("AJ4").Value = 26
$ws.Range("AO4").Value = 65
$ws.Range("U5").Value = 1.64
$ws.Range("X5").Value = 11
$ws.Range("AC5").Value = 11
$ws.Range("AD5").Value = 30
$ws.Range("AF5").Value = 9.800000000000001
$ws.Range("AL5").Value = 70
$ws.Range("F6").Value = 2.46
$ws.Range("Y6").Value = 10.5
$ws.Range("AD6").Value = 18.5
$ws.Range("AF6").Value = 16
$ws.Range("H7").Value = 2.06
$ws.Range("V7").Value = 1.76
$ws.Range("X7").Value = 1000
$ws.Range("Y7").Value = 1000
$ws.Range("Z7").Value = 1000
$ws.Range("AA7").Value = 980
$ws.Range("AB7").Value = 1000
$ws.Range("AC7").Value = 1000
$ws.Range("AD7").Value = 1000
$ws.Range("AE7").Value = 980
$ws.Range("AO7").Value = 980
$ws.Range("F8").Value = 1.96
$ws.Range("G8").Value = 2.22
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = 5
$ws.Range("K8").Value = 3.85
$ws.Range("N8").Value = 3.25
$ws.Range("O8").Value = 1.38
$ws.Range("P8").Value = 1.79
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.26
$ws.Range("T8").Value = 1.71
$ws.Range("U8").Value = 1.85
$ws.Range("V8").Value = 1.25
$ws.Range("W8").Value = 1.82
$ws.Range("Z8").Value = 42
$ws.Range("AF8").Value = 17
$ws.Range("AJ8").Value = 34
$ws.Range("AL8").Value = 980
$ws.Range("AN8").Value = 24
$ws.Range("I9").Value = 6.2
$ws.Range("M9").Value = 1.15
$ws.Range("N9").Value = 2.18
$ws.Range("W9").Value = 1.83
$ws.Range("X9").Value = 1000
$ws.Range("Y9").Value = 1000
$ws.Range("AB9").Value = 1000
$ws.Range("AC9").Value = 1000
$ws.Range("AF9").Value = 1000
$ws.Range("AG9").Value = 1000
$ws.Range("AH9").Value = 980
$ws.Range("AL9").Value = 110
$ws.Range("I10").Value = 1.79
$ws.Range("V10").Value = 2.26
$ws.Range("L11").Value = 1.3
$ws.Range("Z11").Value = 980
$ws.Range("AL11").Value = 980
$ws.Range("T12").Value = 1.78
$ws.Range("F19").Value = 3.65
$ws.Range("K19").Value = 3.9
$ws.Range("P19").Value = 2.34
$ws.Range("V19").Value = 1.87
$ws.Range("X19").Value = 19
$ws.Range("Y19").Value = 12.5
$ws.Range("AE19").Value = 19
$ws.Range("AJ19").Value = 70
$ws.Range("AL19").Value = 42
$ws.Range("AM19").Value = 70
$ws.Range("AN20").Value = 22
$ws.Range("G21").Value = 2.6
$ws.Range("O21").Value = 1.02
$ws.Range("Q21").Value = 2.24
$ws.Range("W21").Value = 1.62
$ws.Range("H22").Value = 3
$ws.Range("K22").Value = 3.35
$ws.Range("Y22").Value = 12
$ws.Range("Z22").Value = 25
$ws.Range("AC22").Value = 8.800000000000001
$ws.Range("AE22").Value = 55
$ws.Range("AH22").Value = 25
$ws.Range("G23").Value = 3.8
$ws.Range("I23").Value = 2.18
$ws.Range("L23").Value = 1.33
$ws.Range("T23").Value = 1.63
$ws.Range("U23").Value = 1.96
$ws.Range("V23").Value = 1.84
$ws.Range("W23").Value = 1.36
$ws.Range("AH23").Value = 21
$ws.Range("AO23").Value = 22
$ws.Range("F24").Value = 3.6
$ws.Range("G24").Value = 3.7
$ws.Range("H24").Value = 2.54
$ws.Range("I24").Value = 2.56
$ws.Range("J24").Value = 3
$ws.Range("K24").Value = 3.05
$ws.Range("AE24").Value = 44
$ws.Range("AM24").Value = 290
